# Natmi following Dr Hou advice
# Recompute LR-pair stats for Efemp1-Egfr to add the "ECs" sending cluster
# (previously only FAPs/sCs were sending clusters), producing a full 3x3
# sending x target cluster matrix and refreshed specificity statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efemp1"
$ws.Range("C2").Value = "Egfr"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.343412
$ws.Range("H2").Value = 1.030236
$ws.Range("I2").Value = 0.004045764821303689
$ws.Range("J2").Value = 0.004045764821303689
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.307106666666667
$ws.Range("N2").Value = 3.92132
$ws.Range("O2").Value = 0.01256263154946851
$ws.Range("P2").Value = 0.01256263154946851
$ws.Range("Q2").Value = 0.4488761146133333
$ws.Range("R2").Value = 4.03988503152
$ws.Range("S2").Value = 0.000050825452785839558703454011
$ws.Range("T2").Value = 0.000050825452785839551927190433

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efemp1"
$ws.Range("C3").Value = "Egfr"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.343412
$ws.Range("H3").Value = 1.030236
$ws.Range("I3").Value = 0.004045764821303689
$ws.Range("J3").Value = 0.004045764821303689
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.22623699999998
$ws.Range("N3").Value = 240.678711
$ws.Range("O3").Value = 0.77105616682495
$ws.Range("P3").Value = 0.77105616682495
$ws.Range("Q3").Value = 27.550652500644
$ws.Range("R3").Value = 247.955872505796
$ws.Range("S3").Value = 0.003119511914989652
$ws.Range("T3").Value = 0.003119511914989651

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efemp1"
$ws.Range("C4").Value = "Egfr"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.343412
$ws.Range("H4").Value = 1.030236
$ws.Range("I4").Value = 0.004045764821303689
$ws.Range("J4").Value = 0.004045764821303689
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 22.51385866666667
$ws.Range("N4").Value = 67.54157600000001
$ws.Range("O4").Value = 0.2163812016255815
$ws.Range("P4").Value = 0.2163812016255815
$ws.Range("Q4").Value = 7.731529232437333
$ws.Range("R4").Value = 69.58376309193601
$ws.Range("S4").Value = 0.0008754274535281982
$ws.Range("T4").Value = 0.000875427453528198

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efemp1"
$ws.Range("C5").Value = "Egfr"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 83.68760400000001
$ws.Range("H5").Value = 251.062812
$ws.Range("I5").Value = 0.9859304981840877
$ws.Range("J5").Value = 0.9859304981840876
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.307106666666667
$ws.Range("N5").Value = 3.92132
$ws.Range("O5").Value = 0.01256263154946851
$ws.Range("P5").Value = 0.01256263154946851
$ws.Range("Q5").Value = 109.38862510576
$ws.Range("R5").Value = 984.4976259518401
$ws.Range("S5").Value = 0.01238588158207063
$ws.Range("T5").Value = 0.01238588158207062

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efemp1"
$ws.Range("C6").Value = "Egfr"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 83.68760400000001
$ws.Range("H6").Value = 251.062812
$ws.Range("I6").Value = 0.9859304981840877
$ws.Range("J6").Value = 0.9859304981840876
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 80.22623699999998
$ws.Range("N6").Value = 240.678711
$ws.Range("O6").Value = 0.77105616682495
$ws.Range("P6").Value = 0.77105616682495
$ws.Range("Q6").Value = 6713.941552466147
$ws.Range("R6").Value = 60425.47397219532
$ws.Range("S6").Value = 0.760207790685636
$ws.Range("T6").Value = 0.7602077906856359

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efemp1"
$ws.Range("C7").Value = "Egfr"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 83.68760400000001
$ws.Range("H7").Value = 251.062812
$ws.Range("I7").Value = 0.9859304981840877
$ws.Range("J7").Value = 0.9859304981840876
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 22.51385866666667
$ws.Range("N7").Value = 67.54157600000001
$ws.Range("O7").Value = 0.2163812016255815
$ws.Range("P7").Value = 0.2163812016255815
$ws.Range("Q7").Value = 1884.130888607968
$ws.Range("R7").Value = 16957.17799747171
$ws.Range("S7").Value = 0.2133368259163811
$ws.Range("T7").Value = 0.2133368259163811

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efemp1"
$ws.Range("C8").Value = "Egfr"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 0.8508333333333332
$ws.Range("H8").Value = 2.5525
$ws.Range("I8").Value = 0.01002373699460868
$ws.Range("J8").Value = 0.01002373699460868
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.307106666666667
$ws.Range("N8").Value = 3.92132
$ws.Range("O8").Value = 0.01256263154946851
$ws.Range("P8").Value = 0.01256263154946851
$ws.Range("Q8").Value = 1.112129922222222
$ws.Range("R8").Value = 10.0091693
$ws.Range("S8").Value = 0.0001259245146120456
$ws.Range("T8").Value = 0.0001259245146120456

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efemp1"
$ws.Range("C9").Value = "Egfr"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 0.8508333333333332
$ws.Range("H9").Value = 2.5525
$ws.Range("I9").Value = 0.01002373699460868
$ws.Range("J9").Value = 0.01002373699460868
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 80.22623699999998
$ws.Range("N9").Value = 240.678711
$ws.Range("O9").Value = 0.77105616682495
$ws.Range("P9").Value = 0.77105616682495
$ws.Range("Q9").Value = 68.25915664749998
$ws.Range("R9").Value = 614.3324098274999
$ws.Range("S9").Value = 0.007728864224324411
$ws.Range("T9").Value = 0.007728864224324411

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efemp1"
$ws.Range("C10").Value = "Egfr"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.8508333333333332
$ws.Range("H10").Value = 2.5525
$ws.Range("I10").Value = 0.01002373699460868
$ws.Range("J10").Value = 0.01002373699460868
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 22.51385866666667
$ws.Range("N10").Value = 67.54157600000001
$ws.Range("O10").Value = 0.2163812016255815
$ws.Range("P10").Value = 0.2163812016255815
$ws.Range("Q10").Value = 19.15554141555555
$ws.Range("R10").Value = 172.39987274
$ws.Range("S10").Value = 0.00216894825567222
$ws.Range("T10").Value = 0.00216894825567222

